$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The date column (A2:A75) needs to shift "up" by one row: each row's date
# becomes the date that used to be one row below it, and a brand new
# (one month older) date is appended at the former bottom (row 75).
# Row 4 also needs to pick up row 5's date number format (yyyy-mm-dd)
# since the "recent" date-format block (rows 2-4) shrinks to just rows 2-3.

$firstRow = 2
$lastRow = 75

# Capture the original values (Excel serial date numbers) before mutating.
$originalValues = @{}
for ($r = $firstRow; $r -le ($lastRow + 1); $r++) {
    $originalValues[$r] = $ws.Cells.Item($r, 1).Value()
}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $originalValues[$r + 1]
}

# New oldest data point, one month before the previous oldest (row 75),
# matching the monthly cadence of the series.
$ws.Cells.Item($lastRow, 1).Value = 43405

# Row 4 now joins the "older" formatting block (rows 5+) instead of the
# "recent" block (rows 2-4). Use the exact escaped format code so the
# existing numFmt/style (the one rows 5+ already use) is reused instead of
# a duplicate number format being minted.
$ws.Cells.Item(4, 1).NumberFormat = "yyyy\-mm\-dd"
